$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-02-04 Sunday"

# Update the division problems in the table, cell by cell
# (direct cell targeting avoids cross-matches since some new values
#  equal other cells old values)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "66÷3=22, 0"
$t.Cell(1, 2).Range.Text = "32÷9=3, 5"
$t.Cell(1, 3).Range.Text = "98÷4=24, 2"
$t.Cell(1, 4).Range.Text = "40÷6=6, 4"
$t.Cell(1, 5).Range.Text = "95÷2=47, 1"
$t.Cell(5, 1).Range.Text = "58÷9=6, 4"
$t.Cell(5, 2).Range.Text = "17÷7=2, 3"
$t.Cell(5, 3).Range.Text = "95÷6=15, 5"
$t.Cell(5, 4).Range.Text = "68÷2=34, 0"
$t.Cell(5, 5).Range.Text = "93÷3=31, 0"
$t.Cell(9, 1).Range.Text = "71÷6=11, 5"
$t.Cell(9, 2).Range.Text = "79÷3=26, 1"
$t.Cell(9, 3).Range.Text = "82÷2=41, 0"
$t.Cell(9, 4).Range.Text = "49÷7=7, 0"
$t.Cell(9, 5).Range.Text = "15÷7=2, 1"
$t.Cell(13, 1).Range.Text = "53÷2=26, 1"
$t.Cell(13, 2).Range.Text = "40÷8=5, 0"
$t.Cell(13, 3).Range.Text = "32÷9=3, 5"
$t.Cell(13, 4).Range.Text = "78÷6=13, 0"
$t.Cell(13, 5).Range.Text = "26÷7=3, 5"
$t.Cell(17, 1).Range.Text = "32÷6=5, 2"
$t.Cell(17, 2).Range.Text = "74÷4=18, 2"
$t.Cell(17, 3).Range.Text = "72÷7=10, 2"
$t.Cell(17, 4).Range.Text = "51÷2=25, 1"
$t.Cell(17, 5).Range.Text = "29÷5=5, 4"
